# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column E) for the row belonging to
# 43759be3-4d94-4262-815f-fcf2bc9b05d8 on both the zh-cn and de-de status
# sheets, reflecting a fresh handoff report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-12 22:34:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-12 22:34:57"
